$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7's date cell (A7) switches from a date-only format to a date+time format
$ws.Range("A7").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new daily row (row 8)
$ws.Range("A8").Value = 45957
$ws.Range("A8").NumberFormat = "YYYY-MM-DD"

$ws.Range("B8").Value = 15
$ws.Range("C8").Value = 21
$ws.Range("D8").Value = 17
